function Set-TextCell {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "E2" "0.97%"
Set-TextCell $ws "D3" "45.56"
Set-TextCell $ws "E3" "3.09%"
Set-TextCell $ws "D4" "5.538"
Set-TextCell $ws "E4" "0.56%"
Set-TextCell $ws "D5" "0.08487"
Set-TextCell $ws "E5" "5.85%"
Set-TextCell $ws "D6" "2.075"
Set-TextCell $ws "E6" "0.65%"
Set-TextCell $ws "D7" "0.9896"
Set-TextCell $ws "E7" "3.66%"
Set-TextCell $ws "E8" "-3.30%"
Set-TextCell $ws "D9" "0.1157"
Set-TextCell $ws "E9" "1.30%"
Set-TextCell $ws "D10" "0.1931"
Set-TextCell $ws "E10" "2.82%"
Set-TextCell $ws "D11" "9.489"
Set-TextCell $ws "E11" "-8.11%"
Set-TextCell $ws "D12" "0.09816"
Set-TextCell $ws "E12" "-0.31%"
Set-TextCell $ws "D13" "0.04692"
Set-TextCell $ws "E13" "-3.55%"
Set-TextCell $ws "E14" "-0.17%"
Set-TextCell $ws "D15" "0.001283"
Set-TextCell $ws "E15" "2.04%"
Set-TextCell $ws "D16" "0.005912"
Set-TextCell $ws "E16" "-2.14%"
Set-TextCell $ws "D17" "3.386"
Set-TextCell $ws "E17" "0.18%"
Set-TextCell $ws "D18" "4.431"
Set-TextCell $ws "E18" "0.68%"
Set-TextCell $ws "E19" "-1.45%"
Set-TextCell $ws "D20" "0.1383"
Set-TextCell $ws "E20" "0.00%"
Set-TextCell $ws "E21" "-1.00%"
Set-TextCell $ws "D22" "0.04142"
Set-TextCell $ws "E22" "1.25%"
Set-TextCell $ws "E23" "-0.08%"
Set-TextCell $ws "D24" "0.004611"
Set-TextCell $ws "E24" "5.75%"
Set-TextCell $ws "D25" "0.0001303"
Set-TextCell $ws "E25" "10.33%"
Set-TextCell $ws "D26" "0.0002985"
Set-TextCell $ws "E26" "-20.32%"
Set-TextCell $ws "D38" "0.02726"
Set-TextCell $ws "E38" "5.87%"
Set-TextCell $ws "D39" "0.05753"
Set-TextCell $ws "E39" "-0.19%"
Set-TextCell $ws "D40" "0.007748"
Set-TextCell $ws "E40" "2.02%"
Set-TextCell $ws "D41" "0.1437"
Set-TextCell $ws "E41" "2.41%"
Set-TextCell $ws "D42" "0.007624"
Set-TextCell $ws "E42" "4.12%"
Set-TextCell $ws "D43" "0.002127"
Set-TextCell $ws "E43" "6.74%"
Set-TextCell $ws "D44" "0.008068"
Set-TextCell $ws "E44" "-10.79%"
Set-TextCell $ws "D45" "0.3556"
Set-TextCell $ws "D46" "0.00007059"
Set-TextCell $ws "E46" "0.80%"
Set-TextCell $ws "E47" "0.20%"
Set-TextCell $ws "E48" "0.22%"
Set-TextCell $ws "B49" "BOLO"
Set-TextCell $ws "C49" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell $ws "D49" "0.003455"
Set-TextCell $ws "E49" "-1.38%"
Set-TextCell $ws "B50" "CoinbaseStockToken"
Set-TextCell $ws "C50" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell $ws "D50" "0.003538"
Set-TextCell $ws "E50" "1.07%"
Set-TextCell $ws "D51" "0.00002105"
Set-TextCell $ws "E51" "0.20%"

Write-Host "Applied 73 cell updates"
